# Apply updated cryptocurrency Price (column D) and Volume(1h) (column E)
# values, matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "30.575.69"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +0.63%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "2.111.65"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +0.18%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "1.013"; ForceText = $true }
    @{ Cell = "E4"; Value = "  +0.90%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "350.03"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +4.81%  "; ForceText = $false }
    @{ Cell = "E6"; Value = "  +0.83%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.5247"; ForceText = $true }
    @{ Cell = "E7"; Value = "  +0.30%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.4507"; ForceText = $true }
    @{ Cell = "D9"; Value = "54.35"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +2.46%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.09009"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +1.01%  "; ForceText = $false }
    @{ Cell = "E11"; Value = "  -0.44%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "24.44"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +0.01%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "2.117.42"; ForceText = $false }
    @{ Cell = "E13"; Value = "  +1.13%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "6.813"; ForceText = $true }
    @{ Cell = "E14"; Value = "  -0.01%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "8.038"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +0.61%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "101.10"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +4.76%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "0.00001168"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +2.91%  "; ForceText = $false }
    @{ Cell = "E18"; Value = "  +0.75%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "0.06715"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +1.31%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "19.39"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +0.64%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "1.010"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +0.84%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "6.287"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -1.27%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "30.650.55"; ForceText = $false }
    @{ Cell = "E23"; Value = "  +0.67%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "12.82"; ForceText = $true }
    @{ Cell = "E24"; Value = "  +3.19%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "2.386"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +1.06%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "2.368.78"; ForceText = $false }
    @{ Cell = "E26"; Value = "  +0.98%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "22.39"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +0.11%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "164.92"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +0.84%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "2.531"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -1.56%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "135.99"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +2.50%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "1.190"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -4.29%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "0.1075"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +0.13%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "1.651"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -4.16%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "6.357"; ForceText = $true }
    @{ Cell = "E34"; Value = "  +0.68%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "4.021"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +2.41%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "10.34"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -1.87%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "5.916"; ForceText = $true }
    @{ Cell = "E37"; Value = "  +6.46%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.02643"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +2.17%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "0.06835"; ForceText = $true }
    @{ Cell = "E39"; Value = "  +0.12%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "0.2306"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +0.39%  "; ForceText = $false }
    @{ Cell = "E41"; Value = "  -1.89%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "0.6869"; ForceText = $true }
    @{ Cell = "E42"; Value = "  -0.58%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "1.270"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +1.83%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "14.65"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +4.75%  "; ForceText = $false }
    @{ Cell = "E45"; Value = "  -0.84%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "0.6443"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +0.96%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "3.757"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +2.74%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "0.00000000359"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +1.78%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "1.251"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +0.34%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "0.07281"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +2.09%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "82.29"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -1.44%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
